# Task3.5 Burndown - update the "Balance (Planned)" burndown figures for
# hours 3, 4, 5 and 9 (rows 5, 6, 7, 11) and leave the cursor on the next
# entry cell (F3), matching the author's re-upload of the tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task 3.5")

$ws.Range("D5").Value = 24
$ws.Range("D6").Value = 21
$ws.Range("D7").Value = 18
$ws.Range("D11").Value = 4

$ws.Range("F3").Select()
